$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.870.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.541.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.00"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.17"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.90%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.539.42"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.113"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.997.47"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.306.18"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.537.81"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.98"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.01"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "329.90"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.24"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.88"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "651.49"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.662.75"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.65%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.89"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.64"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.16"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.84"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.58%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.373"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("B42").Value = "EthereumClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.99"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "162.85"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.39%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.64"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.53"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.624"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.94%  "
